# This script refreshes the "江西-漫展信息" workbook's exhibition listing data,
# mirroring a re-scrape that dropped the 4 oldest (already-passed) events and
# picked up newer "want to go" counts for the remaining, still-listed events.

$wb = $excel.ActiveWorkbook

# New "want to go" (column F) counts, keyed by each event's unique bilibili link
# (column H). These reflect the handful of events whose popularity count moved
# since the previous snapshot.
$fUpdates = @{
    "https://show.bilibili.com/platform/detail.html?id=86684" = 618
    "https://show.bilibili.com/platform/detail.html?id=88294" = 72
    "https://show.bilibili.com/platform/detail.html?id=83980" = 4504
    "https://show.bilibili.com/platform/detail.html?id=87164" = 1809
    "https://show.bilibili.com/platform/detail.html?id=88903" = 124
    "https://show.bilibili.com/platform/detail.html?id=86341" = 3042
    "https://show.bilibili.com/platform/detail.html?id=86683" = 578
    "https://show.bilibili.com/platform/detail.html?id=86453" = 237
    "https://show.bilibili.com/platform/detail.html?id=87449" = 577
    "https://show.bilibili.com/platform/detail.html?id=87225" = 498
    "https://show.bilibili.com/platform/detail.html?id=84407" = 494
    "https://show.bilibili.com/platform/detail.html?id=84102" = 1746
    "https://show.bilibili.com/platform/detail.html?id=84184" = 1283
    "https://show.bilibili.com/platform/detail.html?id=88602" = 110
    "https://show.bilibili.com/platform/detail.html?id=88514" = 1506
    "https://show.bilibili.com/platform/detail.html?id=87600" = 120
    "https://show.bilibili.com/platform/detail.html?id=89411" = 38
    "https://show.bilibili.com/platform/detail.html?id=89821" = 34
    "https://show.bilibili.com/platform/detail.html?id=89466" = 78
    "https://show.bilibili.com/platform/detail.html?id=89295" = 3330
    "https://show.bilibili.com/platform/detail.html?id=87135" = 732
    "https://show.bilibili.com/platform/detail.html?id=89742" = 57
    "https://show.bilibili.com/platform/detail.html?id=89659" = 216
    "https://show.bilibili.com/platform/detail.html?id=89738" = 1631
}

# Both the "展览" (Exhibitions) and "全部类型" (All types) sheets list the same
# events (plus a couple of extra performance-only rows on the latter) and both
# need the same treatment: drop the 4 earliest events (rows 2-5) and refresh
# the "want to go" counts of whichever events remain.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Remove the 4 oldest events; remaining rows shift up automatically.
    $ws.Range("A2:A5").EntireRow.Delete()

    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count

    # Column A holds a manual 0-based sequence number (header row = 0), so it
    # must be renumbered after the shift.
    for ($r = 2; $r -le $rowCount; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }

    # Refresh "want to go" counts (column F) for events whose link we recognize.
    for ($r = 2; $r -le $rowCount; $r++) {
        $link = $ws.Cells.Item($r, 8).Value()
        if ($fUpdates.ContainsKey($link)) {
            $ws.Cells.Item($r, 6).Value = $fUpdates[$link]
        }
    }
}
